# Apply crypto price/volume/name/link updates as described in the commit diff.
# Commit: "Updated cryptos list on Sat May  6 07:09:39 UTC 2023 with GitHub Actions"
#
# All target cells hold text values (prices use "." as thousands separators, e.g.
# "29.521.90", and percentages are padded strings like "  +0.95%  "), so each new
# value is written with a leading apostrophe to force Excel to store it as literal
# text instead of auto-converting it to a number. The apostrophe itself is not part
# of the stored cell value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.521.90"
$ws.Range('E2').Value = "'  +0.95%  "
$ws.Range('D3').Value = "'1.971.43"
$ws.Range('E3').Value = "'  +3.42%  "
$ws.Range('D4').Value = "'1.004"
$ws.Range('E4').Value = "'  +0.34%  "
$ws.Range('D5').Value = "'327.20"
$ws.Range('E5').Value = "'  +0.28%  "
$ws.Range('D6').Value = "'1.004"
$ws.Range('E6').Value = "'  +0.37%  "
$ws.Range('D7').Value = "'0.4657"
$ws.Range('E7').Value = "'  +0.13%  "
$ws.Range('D8').Value = "'0.3907"
$ws.Range('E8').Value = "'  -0.34%  "
$ws.Range('D9').Value = "'46.24"
$ws.Range('E9').Value = "'  +0.32%  "
$ws.Range('D10').Value = "'0.07936"
$ws.Range('E10').Value = "'  +0.48%  "
$ws.Range('D11').Value = "'0.9864"
$ws.Range('E11').Value = "'  -0.59%  "
$ws.Range('E12').Value = "'  +4.08%  "
$ws.Range('D13').Value = "'1.972.81"
$ws.Range('E13').Value = "'  +5.30%  "
$ws.Range('D14').Value = "'7.169"
$ws.Range('E14').Value = "'  +1.20%  "
$ws.Range('D15').Value = "'5.799"
$ws.Range('E15').Value = "'  +0.88%  "
$ws.Range('D16').Value = "'0.07120"
$ws.Range('E16').Value = "'  +1.84%  "
$ws.Range('D17').Value = "'87.72"
$ws.Range('E17').Value = "'  -0.61%  "
$ws.Range('D18').Value = "'1.007"
$ws.Range('E18').Value = "'  +0.51%  "
$ws.Range('D19').Value = "'0.000009906"
$ws.Range('E19').Value = "'  -0.95%  "
$ws.Range('D20').Value = "'17.19"
$ws.Range('E20').Value = "'  +0.33%  "
$ws.Range('D21').Value = "'1.006"
$ws.Range('E21').Value = "'  +0.44%  "
$ws.Range('D22').Value = "'29.532.11"
$ws.Range('E22').Value = "'  +1.00%  "
$ws.Range('B23').Value = "'BitDAO"
$ws.Range('C23').Value = "'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range('D23').Value = "'0.5094"
$ws.Range('E23').Value = "'  +7.01%  "
$ws.Range('B24').Value = "'Uniswap"
$ws.Range('C24').Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range('D24').Value = "'5.537"
$ws.Range('E24').Value = "'  +4.14%  "
$ws.Range('B25').Value = "'Cosmos"
$ws.Range('C25').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('D25').Value = "'11.13"
$ws.Range('E25').Value = "'  +0.29%  "
$ws.Range('B26').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C26').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D26').Value = "'2.225.96"
$ws.Range('E26').Value = "'  +6.18%  "
$ws.Range('B27').Value = "'Toncoin"
$ws.Range('C27').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D27').Value = "'2.109"
$ws.Range('E27').Value = "'  +0.04%  "
$ws.Range('B28').Value = "'Monero"
$ws.Range('C28').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D28').Value = "'158.50"
$ws.Range('E28').Value = "'  +1.45%  "
$ws.Range('B29').Value = "'EthereumClassic"
$ws.Range('C29').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('D29').Value = "'19.51"
$ws.Range('E29').Value = "'  +0.41%  "
$ws.Range('B30').Value = "'InternetComputer(DFINITY)"
$ws.Range('C30').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D30').Value = "'5.771"
$ws.Range('E30').Value = "'  -3.45%  "
$ws.Range('B31').Value = "'BitcoinCash"
$ws.Range('C31').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('D31').Value = "'119.58"
$ws.Range('E31').Value = "'  +0.71%  "
$ws.Range('B32').Value = "'LidoDAOToken"
$ws.Range('C32').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('D32').Value = "'1.881"
$ws.Range('E32').Value = "'  -0.45%  "
$ws.Range('B33').Value = "'Stellar"
$ws.Range('C33').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('D33').Value = "'0.09414"
$ws.Range('E33').Value = "'  +0.61%  "
$ws.Range('B34').Value = "'ImmutableX"
$ws.Range('C34').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D34').Value = "'0.8776"
$ws.Range('E34').Value = "'  -2.87%  "
$ws.Range('B35').Value = "'Filecoin"
$ws.Range('C35').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('D35').Value = "'5.219"
$ws.Range('E35').Value = "'  -0.80%  "
$ws.Range('B36').Value = "'ARBITRUM"
$ws.Range('C36').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D36').Value = "'1.318"
$ws.Range('E36').Value = "'  -0.84%  "
$ws.Range('B37').Value = "'HuobiToken"
$ws.Range('C37').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D37').Value = "'3.134"
$ws.Range('E37').Value = "'  -1.47%  "
$ws.Range('B38').Value = "'Hedera"
$ws.Range('C38').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D38').Value = "'0.05793"
$ws.Range('E38').Value = "'  +0.19%  "
$ws.Range('B39').Value = "'TrustWalletToken"
$ws.Range('C39').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D39').Value = "'1.163"
$ws.Range('E39').Value = "'  -1.67%  "
$ws.Range('B40').Value = "'VeChain"
$ws.Range('C40').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('D40').Value = "'0.02095"
$ws.Range('E40').Value = "'  +0.08%  "
$ws.Range('D41').Value = "'0.5692"
$ws.Range('E41').Value = "'  -0.37%  "
$ws.Range('B42').Value = "'FraxShare"
$ws.Range('C42').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D42').Value = "'7.678"
$ws.Range('E42').Value = "'  -0.46%  "
$ws.Range('B43').Value = "'Algorand"
$ws.Range('C43').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D43').Value = "'0.1792"
$ws.Range('E43').Value = "'  -0.15%  "
$ws.Range('B44').Value = "'Aptos"
$ws.Range('C44').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('D44').Value = "'9.611"
$ws.Range('E44').Value = "'  -1.35%  "
$ws.Range('B45').Value = "'MXToken"
$ws.Range('C45').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D45').Value = "'2.756"
$ws.Range('E45').Value = "'  +7.22%  "
$ws.Range('B46').Value = "'PEPE"
$ws.Range('C46').Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('D46').Value = "'0.000002784"
$ws.Range('E46').Value = "'  +45.83%  "
$ws.Range('B47').Value = "'Decentraland"
$ws.Range('C47').Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range('D47').Value = "'0.5316"
$ws.Range('E47').Value = "'  -0.94%  "
$ws.Range('B48').Value = "'EnergySwap"
$ws.Range('C48').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D48').Value = "'11.64"
$ws.Range('E48').Value = "'  -2.78%  "
$ws.Range('B49').Value = "'RenderToken"
$ws.Range('C49').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D49').Value = "'2.142"
$ws.Range('E49').Value = "'  -1.82%  "
$ws.Range('B50').Value = "'Cronos"
$ws.Range('C50').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D50').Value = "'0.06922"
$ws.Range('E50').Value = "'  -1.37%  "
$ws.Range('B51').Value = "'NEARProtocol"
$ws.Range('C51').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('D51').Value = "'1.824"
$ws.Range('E51').Value = "'  -1.89%  "
